$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "submit pdf or a slide deck" paragraph becomes a green-colored
# instruction to submit a slide deck (powerpoint / impress).
# ---------------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("pdf")) {
        $pSubmit = $p
        $found = $true
        break
    }
}

if ($found) {
    # Paint the whole paragraph (incl. paragraph mark) green first - this
    # updates every existing run plus the pilcrow formatting in one shot.
    $pSubmit.Range.Font.Color = 5287936   # RGB 00B050

    $r1 = $pSubmit.Range
    $r1.Find.Execute("אפשר להגיש קובץ ") | Out-Null
    $r1.Text = "יש להגיש מצגת ("

    $r2 = $pSubmit.Range
    $r2.Find.Execute("pdf") | Out-Null
    $r2.Text = "powerpoint / impress"

    $r3 = $pSubmit.Range
    $r3.Find.Execute(" או מצגת לבחירתכם.") | Out-Null
    $r3.Text = ")."
}

# ---------------------------------------------------------------------------
# Change 2 & 3: the "playtest ... personal experience" paragraph - the
# _GoBack bookmark moves from mid-word (splitting "לכל") to the start of the
# paragraph, and the two runs that spelled out "לכ" + "ל אחד ..." are
# rejoined into a single run.
# ---------------------------------------------------------------------------
$pPlay = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Contains("playtest") -and $t.Contains("חוויית-שחקן שונה")) {
        $pPlay = $p
        break
    }
}

if ($pPlay -ne $null) {
    # Re-anchoring a bookmark with the same name moves it (removes the old
    # occurrence automatically) - point it at the whole paragraph so
    # bookmarkStart lands right before the first run and bookmarkEnd right
    # after the last one.
    $d.Bookmarks.Add("_GoBack", $pPlay.Range) | Out-Null

    $rMerge = $pPlay.Range
    $rMerge.Find.Execute(" באופן אישי – כי לכ") | Out-Null
    $rMerge.Text = " באופן אישי – כי לכל אחד יש חוויית-שחקן שונה."

    $rTail = $pPlay.Range
    $rTail.Find.Execute("ל אחד יש חוויית-שחקן שונה.") | Out-Null
    $rTail.Delete()
}
